# "update misc for hlth raster"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet   # Sheet2 (the only / tab-selected sheet)

# The health-distance raster variable was renamed upstream; fix the
# mislabeled target name stored in the datamap table (row 15 / col A).
$ws.Range("A15").Value = "hlthdist_fctb_clst"

# Re-point the sheet's active selection at the last populated row (A15)
# instead of the empty row below the table (A16).
$ws.Range("A15").Select()

# Remember the window position at save time (as Excel does on close).
$win = $excel.ActiveWindow
$win.Left = 1060
$win.Top = 7780
